# Georgia overview workbook text-edit pass.
# Converts a batch of numeric "count" cells into text cells (so they render
# with thousands separators / match the rest of the text-formatted sheet),
# rewrites the all-zero county rows into percentage/currency text, and adds
# a new "Total" row to the County sheet.

function Set-TextValue($ws, $row, $col, $text) {
    # Force the cell to Text so Excel doesn't silently re-parse a string
    # like "2,130" back into a number, then drop the temporary format so
    # the cell is left with no explicit style (matching the rest of the
    # sheet's plain inline-string cells).
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overall sheet: A2 (990 filer count) becomes text "2,130"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall 2 1 "2,130"

# ---------------------------------------------------------------------
# County sheet: column B counts become text; zero rows become
# percentage/currency text; a new Total row is appended.
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyData = @(
    "2|2", "3|1", "4|7", "5|1", "6|8", "7|12", "8|1", "9|1", "10|69", "11|1",
    "12|1", "13|1", "14|19", "15|1", "16|1", "17|8", "18|16", "19|7", "20|90", "21|2",
    "22|27", "23|44", "24|23", "25|109", "26|3", "27|9", "28|8", "29|1", "30|18", "31|5",
    "32|3", "33|5", "34|220", "35|5", "36|1", "37|35", "38|15", "39|3", "40|5", "41|4",
    "42|7", "43|26", "44|30", "45|28", "46|6", "47|528", "48|2", "49|25", "50|7", "51|3",
    "52|7", "53|109", "54|7", "55|38", "56|1", "57|3", "58|3", "59|1", "60|1", "61|20",
    "62|27", "63|2", "64|4", "65|3", "66|1", "67|1", "68|2", "69|4", "70|2", "71|1",
    "72|6", "73|4", "74|20", "75|9", "76|3", "77|2", "78|1", "79|2", "80|2", "81|3",
    "82|4", "83|2", "84|4", "85|5", "86|58", "87|8", "88|6", "89|10", "90|1", "91|13",
    "92|2", "93|1", "94|3", "95|3", "96|3", "97|5", "98|3", "99|65", "100|18", "101|1",
    "102|2", "103|11", "104|8", "105|1", "106|9", "107|3", "108|22", "109|8", "110|6", "111|5",
    "112|1", "113|21", "114|1", "115|6", "116|2", "117|4", "118|8", "119|10", "120|1", "121|3",
    "122|3", "123|1", "124|7", "125|22", "126|1", "127|1", "128|3", "129|2"
)

foreach ($entry in $countyData) {
    $parts = $entry -split '\|'
    $row = [int]$parts[0]
    Set-TextValue $wsCounty $row 2 $parts[1]
}

# Rows 130-146 are the zero-activity counties: rewrite each of B:F as
# percentage/currency text instead of the bare "0".
for ($row = 130; $row -le 146; $row++) {
    Set-TextValue $wsCounty $row 2 "0.00%"
    Set-TextValue $wsCounty $row 3 "`$0"
    Set-TextValue $wsCounty $row 4 "0.00%"
    Set-TextValue $wsCounty $row 5 "0.00%"
    Set-TextValue $wsCounty $row 6 "0.00%"
}

# New Total row at the bottom of the County sheet.
Set-TextValue $wsCounty 147 1 "Total"
Set-TextValue $wsCounty 147 2 "2,130"
Set-TextValue $wsCounty 147 3 "`$5,229,883,749"
Set-TextValue $wsCounty 147 4 "9.35%"
Set-TextValue $wsCounty 147 5 "-11.79%"
Set-TextValue $wsCounty 147 6 "66.43%"

# ---------------------------------------------------------------------
# Congressional District sheet: column B counts become text.
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")

$cdData = @(
    "2|150", "3|99", "4|119", "5|126", "6|75", "7|100", "8|217", "9|110", "10|168", "11|516",
    "12|123", "13|100", "14|108", "15|119", "16|2,130"
)

foreach ($entry in $cdData) {
    $parts = $entry -split '\|'
    $row = [int]$parts[0]
    Set-TextValue $wsCd $row 2 $parts[1]
}

# ---------------------------------------------------------------------
# Size sheet: column B counts become text.
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

$sizeData = @(
    "2|665", "3|565", "4|395", "5|175", "6|267", "7|63", "8|2,130"
)

foreach ($entry in $sizeData) {
    $parts = $entry -split '\|'
    $row = [int]$parts[0]
    Set-TextValue $wsSize $row 2 $parts[1]
}

# ---------------------------------------------------------------------
# Subsector sheet: column B counts become text.
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")

$subData = @(
    "2|175", "3|277", "4|85", "5|195", "6|36", "7|682", "8|24", "9|3", "10|162", "11|101",
    "12|360", "13|30", "14|2,130"
)

foreach ($entry in $subData) {
    $parts = $entry -split '\|'
    $row = [int]$parts[0]
    Set-TextValue $wsSub $row 2 $parts[1]
}
